$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B
$ws.Range("B2").Value = 13.41359105506008
$ws.Range("B3").Value = 12.68730325286515
$ws.Range("B4").Value = 12.2213536507598
$ws.Range("B5").Value = 12.02664687339946
$ws.Range("B6").Value = 11.99403099497241
$ws.Range("B7").Value = 12.21874702247058
$ws.Range("B8").Value = 13.16742915302498
$ws.Range("B9").Value = 14.86193313754198
$ws.Range("B10").Value = 16.00557053111997
$ws.Range("B11").Value = 16.55283222292811
$ws.Range("B12").Value = 16.75525086228024
$ws.Range("B13").Value = 16.7118702718256
$ws.Range("B14").Value = 16.56958174143848
$ws.Range("B15").Value = 16.48179922527373
$ws.Range("B16").Value = 15.9691302789828
$ws.Range("B17").Value = 15.67649844044708
$ws.Range("B18").Value = 15.50789264595694
$ws.Range("B19").Value = 15.45040980534519
$ws.Range("B20").Value = 15.70751499291993
$ws.Range("B21").Value = 16.61150589014245
$ws.Range("B22").Value = 17.19174447285765
$ws.Range("B23").Value = 16.88462078595725
$ws.Range("B24").Value = 15.69349986023883
$ws.Range("B25").Value = 14.42215365624403

# Column C
$ws.Range("C2").Value = 11.17353494206511
$ws.Range("C3").Value = 10.50790680757608
$ws.Range("C4").Value = 10.07542857226813
$ws.Range("C5").Value = 9.893243026788587
$ws.Range("C6").Value = 9.862633227741521
$ws.Range("C7").Value = 10.07299558317074
$ws.Range("C8").Value = 10.948978132409
$ws.Range("C9").Value = 12.47716869927993
$ws.Range("C10").Value = 13.48365934277577
$ws.Range("C11").Value = 13.91626621163644
$ws.Range("C12").Value = 14.07645452526466
$ws.Range("C13").Value = 14.04211653119122
$ws.Range("C14").Value = 13.92951774698908
$ws.Range("C15").Value = 13.86007501148601
$ws.Range("C16").Value = 13.45487932253697
$ws.Range("C17").Value = 13.19983682533159
$ws.Range("C18").Value = 13.05076613067331
$ws.Range("C19").Value = 12.99988514433148
$ws.Range("C20").Value = 13.22723255120349
$ws.Range("C21").Value = 13.96268925019546
$ws.Range("C22").Value = 14.42219370738832
$ws.Range("C23").Value = 14.17888259020749
$ws.Range("C24").Value = 13.21485455423471
$ws.Range("C25").Value = 12.08414633823331

# Column E
$ws.Range("E2").Value = 16.61294039418812
$ws.Range("E3").Value = 15.66005130527435
$ws.Range("E4").Value = 15.04982640551767
$ws.Range("E5").Value = 14.79509799898178
$ws.Range("E6").Value = 14.75244357188001
$ws.Range("E7").Value = 15.04641518115324
$ws.Range("E8").Value = 16.28974773205969
$ws.Range("E9").Value = 18.64470577457221
$ws.Range("E10").Value = 20.31202941339512
$ws.Range("E11").Value = 21.02828792452709
$ws.Range("E12").Value = 21.29348514200089
$ws.Range("E13").Value = 21.23663798079229
$ws.Range("E14").Value = 21.05022658966803
$ws.Range("E15").Value = 20.93525945222925
$ws.Range("E16").Value = 20.26437385127669
$ws.Range("E17").Value = 19.84202070029046
$ws.Range("E18").Value = 19.59511867536005
$ws.Range("E19").Value = 19.51083796302969
$ws.Range("E20").Value = 19.88739211646664
$ws.Range("E21").Value = 21.10514361662004
$ws.Range("E22").Value = 21.86585855241529
$ws.Range("E23").Value = 21.46305554971618
$ws.Range("E24").Value = 19.86689242814381
$ws.Range("E25").Value = 17.99302267873476

# Column F
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("F5").Value = 15.008197319934
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("F23").Value = 21.82633154475857
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("F25").Value = 18.34778573295691

# Column G
$ws.Range("G2").Value = 24.76782681396966
$ws.Range("G3").Value = 24.21362577844373
$ws.Range("G4").Value = 23.88204491574791
$ws.Range("G5").Value = 23.74935583636965
$ws.Range("G6").Value = 23.72747692576904
$ws.Range("G7").Value = 23.88024525363028
$ws.Range("G8").Value = 24.57508037422891
$ws.Range("G9").Value = 25.99583360599886
$ws.Range("G10").Value = 27.06075736864512
$ws.Range("G11").Value = 27.5468813053567
$ws.Range("G12").Value = 27.73097687056201
$ws.Range("G13").Value = 27.69133165687326
$ws.Range("G14").Value = 27.56202801836584
$ws.Range("G15").Value = 27.4828205853746
$ws.Range("G16").Value = 27.02900187585822
$ws.Range("G17").Value = 26.75086189257885
$ws.Range("G18").Value = 26.59105481944972
$ws.Range("G19").Value = 26.53698362569691
$ws.Range("G20").Value = 26.78045431451113
$ws.Range("G21").Value = 27.60000911430446
$ws.Range("G22").Value = 28.13560140866177
$ws.Range("G23").Value = 27.84982195107749
$ws.Range("G24").Value = 26.76707525342711
$ws.Range("G25").Value = 25.6067572093453

# Column H
$ws.Range("H2").Value = 12.32400480665108
$ws.Range("H3").Value = 12.33782069623732
$ws.Range("H4").Value = 12.35038539427714
$ws.Range("H5").Value = 12.35652588599671
$ws.Range("H6").Value = 12.35760694218565
$ws.Range("H7").Value = 12.35046408465434
$ws.Range("H8").Value = 12.32791788233084
$ws.Range("H9").Value = 12.31633645808192
$ws.Range("H10").Value = 12.32801476608249
$ws.Range("H11").Value = 12.33775958938403
$ws.Range("H12").Value = 12.34208983227784
$ws.Range("H13").Value = 12.3411287334358
$ws.Range("H14").Value = 12.33810299979665
$ws.Range("H15").Value = 12.33633307717945
$ws.Range("H16").Value = 12.32746735892917
$ws.Range("H17").Value = 12.32316592984249
$ws.Range("H18").Value = 12.32110893549742
$ws.Range("H19").Value = 12.32048400971516
$ws.Range("H20").Value = 12.3235806367706
$ws.Range("H21").Value = 12.33897434127075
$ws.Range("H22").Value = 12.35276667652091
$ws.Range("H23").Value = 12.34506329587878
$ws.Range("H24").Value = 12.32339185226984
$ws.Range("H25").Value = 12.31594305105402

# Column I
$ws.Range("I2").Value = 16.74493826921271
$ws.Range("I3").Value = 16.81553792977347
$ws.Range("I4").Value = 16.86520336765764
$ws.Range("I5").Value = 16.88701623821021
$ws.Range("I6").Value = 16.89073286850019
$ws.Range("I7").Value = 16.86549119112882
$ws.Range("I8").Value = 16.76796122638797
$ws.Range("I9").Value = 16.6274680187171
$ws.Range("I10").Value = 16.55604632674722
$ws.Range("I11").Value = 16.53063333552173
$ws.Range("I12").Value = 16.52204128725063
$ws.Range("I13").Value = 16.52384566035787
$ws.Range("I14").Value = 16.52990571775512
$ws.Range("I15").Value = 16.53375237819141
$ws.Range("I16").Value = 16.55785071197146
$ws.Range("I17").Value = 16.57445674084126
$ws.Range("I18").Value = 16.58467343829566
$ws.Range("I19").Value = 16.58824648652764
$ws.Range("I20").Value = 16.57262003948435
$ws.Range("I21").Value = 16.52809763706241
$ws.Range("I22").Value = 16.50501678644097
$ws.Range("I23").Value = 16.51678055323488
$ws.Range("I24").Value = 16.57344832717891
$ws.Range("I25").Value = 16.6599510726421

# Column N
$ws.Range("N2").Value = 15.4909256585513
$ws.Range("N3").Value = 15.55429657686394
$ws.Range("N4").Value = 15.59514497726952
$ws.Range("N5").Value = 15.61228001260857
$ws.Range("N6").Value = 15.61515485134995
$ws.Range("N7").Value = 15.59537408446849
$ws.Range("N8").Value = 15.5123745696632
$ws.Range("N9").Value = 15.36492251163233
$ws.Range("N10").Value = 15.26582434812521
$ws.Range("N11").Value = 15.22272667550269
$ws.Range("N12").Value = 15.20669028197482
$ws.Range("N13").Value = 15.21013140685796
$ws.Range("N14").Value = 15.22140167232325
$ws.Range("N15").Value = 15.22834195016556
$ws.Range("N16").Value = 15.26868066236543
$ws.Range("N17").Value = 15.29393398437594
$ws.Range("N18").Value = 15.30864571492588
$ws.Range("N19").Value = 15.31365895892703
$ws.Range("N20").Value = 15.29122641353905
$ws.Range("N21").Value = 15.21808363147807
$ws.Range("N22").Value = 15.17193396606098
$ws.Range("N23").Value = 15.19641406033126
$ws.Range("N24").Value = 15.29244990419196
$ws.Range("N25").Value = 15.40318345952331
